$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item("Summary")
$wsTest = $wb.Worksheets.Item("Default test")

# --- Update Summary sheet values (Start Time, End Time, Duration) ---
$wsSummary.Range("B6").Value = "2018-07-16T06:31:26Z"
$wsSummary.Range("B7").Value = "2018-07-16T06:31:56Z"
$wsSummary.Range("B8").Value = "30607 ms"

# --- Update Default test sheet values (Exception, Start, End, Duration) ---
$exceptionText = @"
no such element: Unable to locate element: {"method":"css selector","selector":".menu.clearfix > li:nth-child(3) > ul > li:nth-child(2) > ul > li:nth-child(5) > a"}
  (Session info: chrome=67.0.3396.99)
  (Driver info: chromedriver=2.40.565498 (ea082db3280dd6843ebfb08a625e3eb905c4f5ab),platform=Windows NT 10.0.15063 x86_64) (WARNING: The server did not provide any stacktrace information)
Command duration or timeout: 0 milliseconds
For documentation on this error, please visit: http://seleniumhq.org/exceptions/no_such_element.html
Build info: version: '3.11.0', revision: 'e59cfb3', time: '2018-03-11T20:33:08.638Z'
System info: host: 'LIPL-HO-L-088', ip: '192.168.18.116', os.name: 'Windows 10', os.arch: 'amd64', os.version: '10.0', java.version: '1.8.0_171'
Driver info: org.openqa.selenium.chrome.ChromeDriver
Capabilities {acceptInsecureCerts: false, acceptSslCerts: false, applicationCacheEnabled: false, browserConnectionEnabled: false, browserName: chrome, chrome: {chromedriverVersion: 2.40.565498 (ea082db3280dd6..., userDataDir: C:\Users\LOKESH~1.SHA\AppDa...}, cssSelectorsEnabled: true, databaseEnabled: false, handlesAlerts: true, hasTouchScreen: false, javascriptEnabled: true, locationContextEnabled: true, mobileEmulationEnabled: false, nativeEvents: true, networkConnectionEnabled: false, pageLoadStrategy: normal, platform: XP, platformName: XP, rotatable: false, setWindowRect: true, takesHeapSnapshot: true, takesScreenshot: true, unexpectedAlertBehaviour: , unhandledPromptBehavior: , version: 67.0.3396.99, webStorageEnabled: true}
Session ID: d7f08faf953b50c54e2ffda55781eaf5
*** Element info: {Using=css selector, value=.menu.clearfix > li:nth-child(3) > ul > li:nth-child(2) > ul > li:nth-child(5) > a}
"@

$wsTest.Range("C2").Value = $exceptionText
$wsTest.Range("D2").Value = "2018-07-16T06:31:46Z"
$wsTest.Range("E2").Value = "2018-07-16T06:31:56Z"
$wsTest.Range("F2").Value = "9337 ms"

# --- Column width changes ---
$wsSummary.Columns.Item(2).ColumnWidth = 27.46484375

$wsTest.Columns.Item(3).ColumnWidth = 255.0
$wsTest.Columns.Item(4).ColumnWidth = 27.46484375
$wsTest.Columns.Item(5).ColumnWidth = 27.46484375
